$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "37.443.89"
Set-TextValue $ws.Range("E2") "  -0.94%  "

Set-TextValue $ws.Range("D3") "2.053.64"
Set-TextValue $ws.Range("E3") "  -1.26%  "

Set-TextValue $ws.Range("D4") "0.998"
Set-TextValue $ws.Range("E4") "  -0.34%  "

Set-TextValue $ws.Range("D5") "231.03"
Set-TextValue $ws.Range("E5") "  -1.04%  "

Set-TextValue $ws.Range("E6") "  -1.59%  "

Set-TextValue $ws.Range("E7") "  +0.08%  "

Set-TextValue $ws.Range("D8") "57.21"
Set-TextValue $ws.Range("E8") "  -1.98%  "

Set-TextValue $ws.Range("D9") "0.387"
Set-TextValue $ws.Range("E9") "  -1.86%  "

Set-TextValue $ws.Range("E10") "  +3.29%  "

Set-TextValue $ws.Range("E11") "  -2.05%  "

Set-TextValue $ws.Range("D12") "14.72"
Set-TextValue $ws.Range("E12") "  -2.20%  "

Set-TextValue $ws.Range("D13") "2.358.27"
Set-TextValue $ws.Range("E13") "  -1.19%  "

Set-TextValue $ws.Range("E14") "  -1.94%  "

Set-TextValue $ws.Range("D15") "0.761"
Set-TextValue $ws.Range("E15") "  -2.29%  "

Set-TextValue $ws.Range("D16") "5.32"
Set-TextValue $ws.Range("E16") "  -1.01%  "

Set-TextValue $ws.Range("D17") "2.064.68"
Set-TextValue $ws.Range("E17") "  -1.21%  "

Set-TextValue $ws.Range("D18") "37.334.29"
Set-TextValue $ws.Range("E18") "  -1.07%  "

Set-TextValue $ws.Range("D19") "6.07"
Set-TextValue $ws.Range("E19") "  -1.03%  "

Set-TextValue $ws.Range("D20") "70.00"
Set-TextValue $ws.Range("E20") "  -1.73%  "

Set-TextValue $ws.Range("D21") "0.0₃0841"
Set-TextValue $ws.Range("E21") "  +0.37%  "

Set-TextValue $ws.Range("D22") "227.40"
Set-TextValue $ws.Range("E22") "  -1.16%  "

Set-TextValue $ws.Range("E23") "  +0.17%  "

Set-TextValue $ws.Range("D24") "2.37"
Set-TextValue $ws.Range("E24") "  -1.14%  "

Set-TextValue $ws.Range("E25") "  -4.02%  "

Set-TextValue $ws.Range("D26") "9.58"
Set-TextValue $ws.Range("E26") "  -2.41%  "

Set-TextValue $ws.Range("D27") "167.80"
Set-TextValue $ws.Range("E27") "  -2.44%  "

Set-TextValue $ws.Range("D28") "1.41"
Set-TextValue $ws.Range("E28") "  +0.47%  "

Set-TextValue $ws.Range("D29") "0.129"
Set-TextValue $ws.Range("E29") "  -4.75%  "

Set-TextValue $ws.Range("E30") "  -2.47%  "

Set-TextValue $ws.Range("D31") "0.118"
Set-TextValue $ws.Range("E31") "  -2.46%  "

Set-TextValue $ws.Range("E32") "  -3.34%  "

Set-TextValue $ws.Range("D33") "4.61"
Set-TextValue $ws.Range("E33") "  -1.50%  "

Set-TextValue $ws.Range("D34") "0.0616"
Set-TextValue $ws.Range("E34") "  -2.67%  "

Set-TextValue $ws.Range("D35") "2.42"
Set-TextValue $ws.Range("E35") "  -2.03%  "

Set-TextValue $ws.Range("E36") "  +0.05%  "

Set-TextValue $ws.Range("E37") "  +0.08%  "

Set-TextValue $ws.Range("D38") "3.25"
Set-TextValue $ws.Range("E38") "  -4.33%  "

Set-TextValue $ws.Range("D39") "5.41"
Set-TextValue $ws.Range("E39") "  -1.02%  "

Set-TextValue $ws.Range("D40") "0.0222"
Set-TextValue $ws.Range("E40") "  -5.57%  "

Set-TextValue $ws.Range("D41") "17.12"
Set-TextValue $ws.Range("E41") "  +1.70%  "

Set-TextValue $ws.Range("D42") "1.490.45"
Set-TextValue $ws.Range("E42") "  +2.53%  "

Set-TextValue $ws.Range("E43") "  -1.20%  "

Set-TextValue $ws.Range("D44") "0.0946"
Set-TextValue $ws.Range("E44") "  -2.73%  "

Set-TextValue $ws.Range("D45") "96.84"
Set-TextValue $ws.Range("E45") "  -5.63%  "

Set-TextValue $ws.Range("E46") "  +1.44%  "

Set-TextValue $ws.Range("D47") "1.02"
Set-TextValue $ws.Range("E47") "  -3.72%  "

Set-TextValue $ws.Range("D48") "7.17"
Set-TextValue $ws.Range("E48") "  -2.20%  "

Set-TextValue $ws.Range("E49") "  -2.25%  "

# Row 50/51: RocketPoolETH and FTXToken swap positions with new values
Set-TextValue $ws.Range("B50") "FTXToken"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D50") "3.74"
Set-TextValue $ws.Range("E50") "  -8.84%  "

Set-TextValue $ws.Range("B51") "RocketPoolETH"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D51") "2.244.44"
Set-TextValue $ws.Range("E51") "  -1.20%  "
